# Apply the edits described by the commit:
#  - Sheet1 renamed to Data
#  - The query-table-backed "M_Code_Method" table (E10:H15, fed by the
#    'M Code Method' Power Query) is removed from the sheet, along with its
#    hidden ExternalData_1 defined name
#  - Selection moved to H11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hidden defined name that pointed at the query table's range.
foreach ($n in $wb.Names) {
    $n.Delete()
}

# Remove the query-result table; this also clears the E10:H15 cells that
# held the refreshed query output.
$queryTable = $ws.ListObjects.Item("M_Code_Method")
$queryTable.Delete()

# Rename the worksheet.
$ws.Name = "Data"

# Move the active selection.
$ws.Range("H11").Select()
